$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '27.800.40'
Set-TextValue 'E2' '  +2.87%  '
Set-TextValue 'D3' '1.866.67'
Set-TextValue 'E3' '  +2.54%  '
Set-TextValue 'D4' '1.040'
Set-TextValue 'E4' '  +2.93%  '
Set-TextValue 'D5' '324.56'
Set-TextValue 'E5' '  +3.12%  '
Set-TextValue 'D6' '1.036'
Set-TextValue 'E6' '  +2.66%  '
Set-TextValue 'D7' '0.4424'
Set-TextValue 'E7' '  +2.58%  '
Set-TextValue 'D8' '0.3801'
Set-TextValue 'E8' '  +2.74%  '
Set-TextValue 'E9' '  +2.71%  '
Set-TextValue 'D10' '0.8856'
Set-TextValue 'E10' '  +1.74%  '
Set-TextValue 'D11' '21.78'
Set-TextValue 'E11' '  +2.03%  '
Set-TextValue 'D12' '1.877.70'
Set-TextValue 'E12' '  -12.52%  '
Set-TextValue 'D13' '5.563'
Set-TextValue 'E14' '  +1.65%  '
Set-TextValue 'D15' '0.07216'
Set-TextValue 'E15' '  +3.55%  '
Set-TextValue 'D16' '83.83'
Set-TextValue 'E16' '  +3.21%  '
Set-TextValue 'E17' '  +3.10%  '
Set-TextValue 'D18' '0.000009156'
Set-TextValue 'E18' '  +3.20%  '
Set-TextValue 'E19' '  +2.67%  '
Set-TextValue 'D20' '15.57'
Set-TextValue 'E20' '  +1.35%  '
Set-TextValue 'D21' '27.827.17'
Set-TextValue 'E21' '  +2.74%  '
Set-TextValue 'E22' '  +2.06%  '
Set-TextValue 'D23' '11.38'
Set-TextValue 'E23' '  +3.17%  '
Set-TextValue 'E24' '  +4.01%  '
Set-TextValue 'D25' '158.60'
Set-TextValue 'E25' '  +2.58%  '
Set-TextValue 'D26' '18.90'
Set-TextValue 'E26' '  +2.40%  '
Set-TextValue 'D27' '1.990'
Set-TextValue 'E27' '  +3.39%  '
Set-TextValue 'D28' '5.323'
Set-TextValue 'E28' '  +1.29%  '
Set-TextValue 'D29' '117.75'
Set-TextValue 'E29' '  +2.36%  '
Set-TextValue 'D30' '0.09111'
Set-TextValue 'E30' '  +1.33%  '
Set-TextValue 'B31' 'ImmutableX'
Set-TextValue 'C31' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D31' '0.7807'
Set-TextValue 'E31' '  +4.28%  '
Set-TextValue 'B32' 'ARBITRUM'
Set-TextValue 'C32' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D32' '1.224'
Set-TextValue 'E32' '  +3.94%  '
Set-TextValue 'D33' '3.058'
Set-TextValue 'E33' '  +8.33%  '
Set-TextValue 'D34' '4.590'
Set-TextValue 'E34' '  +3.48%  '
Set-TextValue 'D35' '1.038'
Set-TextValue 'E35' '  +2.82%  '
Set-TextValue 'D36' '1.163'
Set-TextValue 'E36' '  +3.29%  '
Set-TextValue 'D37' '0.02002'
Set-TextValue 'E37' '  +3.69%  '
Set-TextValue 'E38' '  +2.22%  '
Set-TextValue 'B39' 'MXToken'
Set-TextValue 'C39' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D39' '2.855'
Set-TextValue 'E39' '  +3.59%  '
Set-TextValue 'B40' 'TheSandbox'
Set-TextValue 'C40' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D40' '0.5208'
Set-TextValue 'E40' '  +1.67%  '
Set-TextValue 'E41' '  +2.49%  '
Set-TextValue 'D42' '6.906'
Set-TextValue 'E42' '  +5.96%  '
Set-TextValue 'D43' '8.719'
Set-TextValue 'E43' '  +4.60%  '
Set-TextValue 'D44' '109.80'
Set-TextValue 'E44' '  +2.23%  '
Set-TextValue 'D45' '10.73'
Set-TextValue 'E45' '  +2.28%  '
Set-TextValue 'D46' '1.727'
Set-TextValue 'E46' '  +4.13%  '
Set-TextValue 'D47' '0.4721'
Set-TextValue 'E47' '  +3.14%  '
Set-TextValue 'D48' '0.06456'
Set-TextValue 'E48' '  +3.52%  '
Set-TextValue 'D49' '1.895'
Set-TextValue 'E49' '  +3.10%  '
Set-TextValue 'D51' '64.65'
Set-TextValue 'E51' '  +2.17%  '
